$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text: "Nombre" -> "name"
$ws.Range("A1").Value = "name"

# Remove participants 23-35 (clear their cell contents, shrinking the
# sheet's used range from A1:A35 down to A1:A22)
$ws.Range("A23:A35").ClearContents()

# Reflect the new selection / scroll position left behind in the sheet
# view after that edit (user had scrolled down to row 16 and selected
# the now-empty A23:A36 block)
$ws.Range("A23:A36").Select()
$excel.ActiveWindow.ScrollRow = 16
